$wb = $excel.ActiveWorkbook

# Duplicate "Bus_Makhulu_f" to create the new multi-axle truck sheet, placed right after it.
$srcSheet = $wb.Worksheets.Item("Bus_Makhulu_f")
$srcSheet.Copy($null, $srcSheet)

# The copy is placed immediately after the source sheet.
$newSheet = $wb.Worksheets.Item($srcSheet.Index + 1)
$newSheet.Name = "Truck_Amandla_A1"

# Update the instance label (H3) for the new vehicle instance.
$newSheet.Range("H3").Value = "Ackermann_Amandla_A1"

# Update the Ackermann geometry parameters for the new truck instance.
$newSheet.Range("F6").Value = -0.97309999999999997
$newSheet.Range("G6").Value = 0.55801000000000001
$newSheet.Range("H6").Value = 2.5924

# Make the newly added sheet the active tab, as it is the one just edited.
$newSheet.Activate()
[void]$newSheet.Range("J11").Select()
